$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update data values
$ws.Range("G3").Value = 30641659270
$ws.Range("M3").Value = 305

# Update selection / view state
$ws.Activate()
$ws.Range("A3").Select()
